$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column G ("The") - shifts H1:Z1 labels left by one (to G1:Y1),
# and drops one trailing 0 per data row (G:Z -> G:Y)
$ws.Columns("G:G").Delete()

# Update Title / Authors / Year / DOI / Access Type for rows 2-11
$ws.Range("B2").Value = "Cyber-Flirting: Playing at Love on the Internet"
$ws.Range("C2").Value = "Monica Therese Whitty"
$ws.Range("D2").Value = "'2003"
$ws.Range("E2").Value = "10.1177/0959354303013003003"
$ws.Range("F2").Value = "Restricted"

$ws.Range("B3").Value = "Cyber scares and prophylactic policies: Crossnational evidence on the effect of cyberattacks on public support for surveillance"
$ws.Range("C3").Value = "Amelia C Arsenault, Sarah E Kreps, Keren LG Snider, Daphna Canetti"
$ws.Range("D3").Value = "'2024"
$ws.Range("E3").Value = "10.1177/00223433241233960"
$ws.Range("F3").Value = "Restricted"

$ws.Range("B4").Value = "Indonesia’s Handling of Terrorists’ Cyber Activities: How Repressive Measures Still Fall Short"
$ws.Range("C4").Value = "Ali Abdullah Wibisono, Rachel Kumendong, Iwa Maulana"
$ws.Range("D4").Value = "'2025"
$ws.Range("E4").Value = "10.1177/23477970241298764"
$ws.Range("F4").Value = "Restricted"

$ws.Range("B5").Value = "Mapping Global Cyberterror Networks: An Empirical Study of Al-Qaeda and ISIS Cyberterrorism Events"
$ws.Range("C5").Value = "Claire Seungeun Lee, Kyung-Shick Choi, Ryan Shandler, Chris Kayser"
$ws.Range("D5").Value = "'2021"
$ws.Range("E5").Value = "10.1177/10439862211001606"
$ws.Range("F5").Value = "Restricted"

$ws.Range("B6").Value = "Tech titans, cyber commons and the war in Ukraine: An incipient shift in international relations"
$ws.Range("C6").Value = "Eviatar Matania, Udi Sommer"
$ws.Range("D6").Value = "'2023"
$ws.Range("E6").Value = "10.1177/00471178231211500"
$ws.Range("F6").Value = "Open Access"

$ws.Range("B7").Value = "How the process of discovering cyberattacks biases our understanding of cybersecurity"
$ws.Range("C7").Value = "Harry Oppenheimer"
$ws.Range("D7").Value = "'2024"
$ws.Range("E7").Value = "10.1177/00223433231217687"
$ws.Range("F7").Value = "Open Access"

$ws.Range("B8").Value = "Cyclones in cyberspace: Information shaping and denial in the 2008 Russia–Georgia war"
$ws.Range("C8").Value = "Ronald J. Deibert, Rafal Rohozinski, Masashi Crete-Nishihata"
$ws.Range("D8").Value = "'2012"
$ws.Range("E8").Value = "10.1177/0967010611431079"
$ws.Range("F8").Value = "Restricted"

$ws.Range("B9").Value = "Wargaming the use of intermediate force capabilities in the gray zone"
$ws.Range("C9").Value = "Kyle D Christensen, Peter Dobias"
$ws.Range("D9").Value = "'2024"
$ws.Range("E9").Value = "10.1177/15485129211010227"
$ws.Range("F9").Value = "Restricted"

$ws.Range("B10").Value = "Cyberattacks and public opinion – The effect of uncertainty in guiding preferences"
$ws.Range("C10").Value = "Eric Jardine, Nathaniel Porter, Ryan Shandler"
$ws.Range("D10").Value = "'2024"
$ws.Range("E10").Value = "10.1177/00223433231218178"
$ws.Range("F10").Value = "Restricted"

$ws.Range("B11").Value = "Global versus Local Optimization in Redundancy Resolution of Robotic Manipulators"
$ws.Range("C11").Value = "Kazem Kazerounian, Zhaoyu Wang"
$ws.Range("D11").Value = "'1988"
$ws.Range("E11").Value = "10.1177/027836498800700501"
$ws.Range("F11").Value = "Restricted"
